$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.331.83"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.666.57"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5339"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2661"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06393"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07864"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.561"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "1.665.82"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5538"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "0.0₅8187"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "26.358.92"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.679"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.040"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1228"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.207"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.500"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05847"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.283"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.590"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.281"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9706"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.829"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.420"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5830"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01607"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8620"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("D41").Value = "1.065.21"
$ws.Range("E41").Value = "  +3.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.833"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "1.806.19"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.015"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.31%  "
$ws.Range("E48").Value = "  -6.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4394"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.024"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05165"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.47%  "
